# Weekly update: insert the latest "Betarraga" price record for
# Macroferia Regional de Talca as a new row 666, pushing the existing
# historical rows (666..702) down to (667..703).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 666 (shifts rows 666:702 -> 667:703).
$ws.Rows.Item(666).Insert()

# Populate the new row with this week's record.
$ws.Cells.Item(666, 1).Value2 = 5
$ws.Cells.Item(666, 2).Value2 = "Macroferia Regional de Talca"
$ws.Cells.Item(666, 3).Value2 = "Maule"
$ws.Cells.Item(666, 4).Value2 = 45267
$ws.Cells.Item(666, 5).Value2 = 7
$ws.Cells.Item(666, 6).Value2 = 100114014
$ws.Cells.Item(666, 7).Value2 = "Betarraga"
$ws.Cells.Item(666, 8).Value2 = "Sin especificar"
$ws.Cells.Item(666, 9).Value2 = "Primera"
$ws.Cells.Item(666, 10).Value2 = 5000
$ws.Cells.Item(666, 11).Value2 = 700
$ws.Cells.Item(666, 12).Value2 = 700
$ws.Cells.Item(666, 13).Value2 = 700
$ws.Cells.Item(666, 14).Value2 = "`$/paquete 5 unidades"
$ws.Cells.Item(666, 15).Value2 = "Región del Maule"
$ws.Cells.Item(666, 16).Value2 = 140
$ws.Cells.Item(666, 17).Value2 = 5
$ws.Cells.Item(666, 18).Value2 = "Hortaliza"
